$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '26.028.11'
Set-TextValue 'E2' '  +0.35%  '
Set-TextValue 'D3' '1.643.41'
Set-TextValue 'E3' '  +0.63%  '
Set-TextValue 'E4' '  +0.73%  '
Set-TextValue 'D5' '216.49'
Set-TextValue 'E5' '  +0.85%  '
Set-TextValue 'E6' '  +0.86%  '
Set-TextValue 'E7' '  +0.73%  '
Set-TextValue 'E8' '  +0.37%  '
Set-TextValue 'D9' '0.0639'
Set-TextValue 'E9' '  +1.09%  '
Set-TextValue 'D10' '19.64'
Set-TextValue 'E10' '  -0.27%  '
Set-TextValue 'E11' '  +0.79%  '
Set-TextValue 'B12' 'WrappedEther'
Set-TextValue 'C12' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D12' '1.691.46'
Set-TextValue 'E12' '  +3.88%  '
Set-TextValue 'B13' 'WrappedliquidstakedEther2.0'
Set-TextValue 'C13' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D13' '1.873.09'
Set-TextValue 'E13' '  +0.73%  '
Set-TextValue 'B14' 'Polkadot'
Set-TextValue 'C14' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D14' '4.29'
Set-TextValue 'E14' '  +1.23%  '
Set-TextValue 'E15' '  -0.05%  '
Set-TextValue 'D16' '0.0₃0765'
Set-TextValue 'E16' '  +1.19%  '
Set-TextValue 'D17' '63.30'
Set-TextValue 'E17' '  +0.69%  '
Set-TextValue 'D18' '26.119.85'
Set-TextValue 'E18' '  +0.75%  '
Set-TextValue 'E19' '  +0.73%  '
Set-TextValue 'D20' '193.22'
Set-TextValue 'E20' '  -0.16%  '
Set-TextValue 'D21' '4.34'
Set-TextValue 'E21' '  -0.96%  '
Set-TextValue 'D22' '9.93'
Set-TextValue 'E22' '  -0.43%  '
Set-TextValue 'D23' '6.23'
Set-TextValue 'E23' '  -0.62%  '
Set-TextValue 'E24' '  +0.67%  '
Set-TextValue 'B25' 'BinanceUSD'
Set-TextValue 'C25' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 'D25' '1.01'
Set-TextValue 'E25' '  +1.13%  '
Set-TextValue 'B26' 'Stellar'
Set-TextValue 'C26' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D26' '0.131'
Set-TextValue 'E26' '  +4.12%  '
Set-TextValue 'D27' '144.13'
Set-TextValue 'E27' '  +1.15%  '
Set-TextValue 'E28' '  +0.43%  '
Set-TextValue 'D29' '15.54'
Set-TextValue 'E29' '  +0.36%  '
Set-TextValue 'E30' '  +1.15%  '
Set-TextValue 'E31' '  -0.59%  '
Set-TextValue 'D32' '3.26'
Set-TextValue 'E32' '  +1.21%  '
Set-TextValue 'D33' '3.28'
Set-TextValue 'E33' '  -0.84%  '
Set-TextValue 'E34' '  -3.47%  '
Set-TextValue 'E35' '  +1.46%  '
Set-TextValue 'D36' '0.904'
Set-TextValue 'E36' '  +0.15%  '
Set-TextValue 'D37' '1.132.18'
Set-TextValue 'E37' '  -0.37%  '
Set-TextValue 'E38' '  -1.97%  '
Set-TextValue 'E39' '  -0.48%  '
Set-TextValue 'E40' '  +0.17%  '
Set-TextValue 'D41' '5.49'
Set-TextValue 'E41' '  +0.35%  '
Set-TextValue 'D42' '99.55'
Set-TextValue 'E42' '  +0.42%  '
Set-TextValue 'E43' '  -0.75%  '
Set-TextValue 'D44' '1.782.89'
Set-TextValue 'E44' '  +0.82%  '
Set-TextValue 'E45' '  +5.08%  '
Set-TextValue 'D46' '56.66'
Set-TextValue 'E46' '  +0.81%  '
Set-TextValue 'D47' '0.0529'
Set-TextValue 'E47' '  +0.56%  '
Set-TextValue 'D48' '1.45'
Set-TextValue 'E48' '  -0.42%  '
Set-TextValue 'D49' '7.71'
Set-TextValue 'E49' '  +1.06%  '
Set-TextValue 'E50' '  +0.35%  '
Set-TextValue 'E51' '  -0.60%  '

Write-Output "Applied cryptos list update."
